$wb = $excel.ActiveWorkbook
$wsPlayable = $wb.Worksheets.Item("Playable (untested)")
$wsAll = $wb.Worksheets.Item("ALL")

# The "ambush" entry (row 1616 on the ALL sheet) is confirmed playable, so it
# is copied onto the "Playable (untested)" sheet as new row 294 ...
$wsPlayable.Cells.Item(294, 1).Value = 294
$wsPlayable.Cells.Item(294, 2).Value = "ambush"
$wsPlayable.Cells.Item(294, 3).Value = "ambush.c"
$wsPlayable.Cells.Item(294, 4).Value = "Z80"
$wsPlayable.Cells.Item(294, 8).Value = "2xAY-8910"
$wsPlayable.Cells.Item(294, 13).Value = "Ambush"

# ... and removed from the ALL sheet (row 1616), shifting every following
# row up by one.
[void]$wsAll.Rows.Item(1616).Delete()

# Rebuild the AutoFilter over the new (smaller) used range.
$wsAll.AutoFilterMode = $false
[void]$wsAll.Range("A1:M1751").AutoFilter()

# Keep the workbook-level defined names in sync with the new last row.
foreach ($n in $wb.Names) {
    $nm = $n.Name
    if ($nm -eq "ALL!_FilterDatabase") {
        $n.RefersTo = "=ALL!`$A`$1:`$M`$1751"
    }
    if ($nm -eq "ALL!LIST") {
        $n.RefersTo = "=ALL!`$B`$1:`$M`$1751"
    }
}

# Update on-screen selections to match where the edit happened.
[void]$wsPlayable.Range("A291:A294").Select()
[void]$wsAll.Range("A1616:XFD1616").Select()

# The ALL sheet was the active tab when the workbook was saved.
[void]$wsAll.Activate()
